$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 0.020954251304430627
$ws.Cells.Item(1, 2).Value = 0.0019334743024129544
$ws.Cells.Item(1, 3).Value = 0.27118751496548488
$ws.Cells.Item(1, 4).Value = 0.035447433579969059
$ws.Cells.Item(1, 5).Value = 0.0000000000000000000095000000014450517
$ws.Cells.Item(1, 6).Value = 0.07088503100484006
$ws.Cells.Item(2, 1).Value = 0.020792715039745994
$ws.Cells.Item(2, 2).Value = 0.0075525255228047709
$ws.Cells.Item(2, 3).Value = 0.34458786365878963
$ws.Cells.Item(2, 4).Value = 0.043436627534130766
$ws.Cells.Item(2, 5).Value = 0.000000000000000000008749292871397808
$ws.Cells.Item(2, 6).Value = 0.70148236344073356
$ws.Cells.Item(3, 1).Value = 0.020776753854910148
$ws.Cells.Item(3, 2).Value = 0.0036908176822399641
$ws.Cells.Item(3, 3).Value = 0.3154407918743854
$ws.Cells.Item(3, 4).Value = 0.036961989125356136
$ws.Cells.Item(3, 5).Value = 0.0000000000000000000041968182960211297
$ws.Cells.Item(3, 6).Value = 0.2199337113034861
$ws.Cells.Item(4, 1).Value = 0.02060640479520217
$ws.Cells.Item(4, 2).Value = 0.099999656379555546
$ws.Cells.Item(4, 3).Value = 0.29159641137925418
$ws.Cells.Item(4, 4).Value = 0.081307693025371885
$ws.Cells.Item(4, 5).Value = 0.0000000000000000000095000000014450517
$ws.Cells.Item(4, 6).Value = 6.3413747062704147
$ws.Cells.Item(5, 1).Value = 0.020644654387038813
$ws.Cells.Item(5, 2).Value = 0.065532313566698058
$ws.Cells.Item(5, 3).Value = 0.29997115767752897
$ws.Cells.Item(5, 4).Value = 0.099999517245935476
$ws.Cells.Item(5, 5).Value = 0.0000000000000000000087535711196777548
$ws.Cells.Item(5, 6).Value = 7.1432888509031276
$ws.Cells.Item(6, 1).Value = 0.020724444369064007
$ws.Cells.Item(6, 2).Value = 0.099999937569161021
$ws.Cells.Item(6, 3).Value = 0.25647629861948207
$ws.Cells.Item(6, 4).Value = 0.076785104936291468
$ws.Cells.Item(6, 5).Value = 0.0000000000000000000072636696250072512
$ws.Cells.Item(6, 6).Value = 4.1844737157476963
$ws.Cells.Item(7, 1).Value = 0.020802633550229745
$ws.Cells.Item(7, 2).Value = 0.010333579679974697
$ws.Cells.Item(7, 3).Value = 0.24852010449952336
$ws.Cells.Item(7, 4).Value = 0.040137971161844969
$ws.Cells.Item(7, 5).Value = 0.0000000000000000000095000000014450517
$ws.Cells.Item(7, 6).Value = 0.085551264334523758
$ws.Cells.Item(8, 1).Value = 0.020774086950160864
$ws.Cells.Item(8, 2).Value = 0.018755405988790851
$ws.Cells.Item(8, 3).Value = 0.29525379105628424
$ws.Cells.Item(8, 4).Value = 0.042391227171831587
$ws.Cells.Item(8, 5).Value = 0.000000000000000000009042963619420154
$ws.Cells.Item(8, 6).Value = 0.27756887113401468
$ws.Cells.Item(9, 1).Value = 0.020583333615276356
$ws.Cells.Item(9, 2).Value = 0.075098718840664533
$ws.Cells.Item(9, 3).Value = 0.28350538250177254
$ws.Cells.Item(9, 4).Value = 0.099999999999975372
$ws.Cells.Item(9, 5).Value = 0.0000000000000000000095000000014450517
$ws.Cells.Item(9, 6).Value = 9.9999999999999787
$ws.Cells.Item(10, 1).Value = 0.020493872021246078
$ws.Cells.Item(10, 2).Value = 0.07719453401994987
$ws.Cells.Item(10, 3).Value = 0.14845539420091575
$ws.Cells.Item(10, 4).Value = 0.099999999999977801
$ws.Cells.Item(10, 5).Value = 0.0000000000000000000095470198600632456
$ws.Cells.Item(10, 6).Value = 9.9999999999999787
$ws.Cells.Item(11, 1).Value = 0.020834701831796629
$ws.Cells.Item(11, 2).Value = 0.010151936798452567
$ws.Cells.Item(11, 3).Value = 0.35339942366537686
$ws.Cells.Item(11, 4).Value = 0.040356746056164246
$ws.Cells.Item(11, 5).Value = 0.0000000000000000000095000000014450517
$ws.Cells.Item(11, 6).Value = 0.099915430088768867
$ws.Cells.Item(12, 1).Value = 0.020511488670741804
$ws.Cells.Item(12, 2).Value = 0.099999999619139954
$ws.Cells.Item(12, 3).Value = 0.28690292163137665
$ws.Cells.Item(12, 4).Value = 0.060436900871091569
$ws.Cells.Item(12, 5).Value = 0.0000000000000000000095000000014450517
$ws.Cells.Item(12, 6).Value = 2.6117377019627193
$ws.Cells.Item(13, 1).Value = 0.020881092545712021
$ws.Cells.Item(13, 2).Value = 0.010484976264411552
$ws.Cells.Item(13, 3).Value = 0.34116996419734591
$ws.Cells.Item(13, 4).Value = 0.041259882090867894
$ws.Cells.Item(13, 5).Value = 0.000000000000000000007978865982613023
$ws.Cells.Item(13, 6).Value = 0.26575970155722622
$ws.Cells.Item(14, 1).Value = 0.02084143666683962
$ws.Cells.Item(14, 2).Value = 0.004879726997613164
$ws.Cells.Item(14, 3).Value = 0.29001148546021854
$ws.Cells.Item(14, 4).Value = 0.035567381956172796
$ws.Cells.Item(14, 5).Value = 0.0000000000000000000095000000014450517
$ws.Cells.Item(14, 6).Value = 0.049112729606687631
$ws.Cells.Item(15, 1).Value = 0.020795464228114666
$ws.Cells.Item(15, 2).Value = 0.051902398367279244
$ws.Cells.Item(15, 3).Value = 0.37936664060168895
$ws.Cells.Item(15, 4).Value = 0.050443731448181574
$ws.Cells.Item(15, 5).Value = 0.0000000000000000000095000103986317825
$ws.Cells.Item(15, 6).Value = 1.567955757161916
$ws.Cells.Item(16, 1).Value = 0.020890828346247479
$ws.Cells.Item(16, 2).Value = 0.02129167751441376
$ws.Cells.Item(16, 3).Value = 0.34459737677321295
$ws.Cells.Item(16, 4).Value = 0.059143731707169779
$ws.Cells.Item(16, 5).Value = 0.0000000000000000000035668904311534509
$ws.Cells.Item(16, 6).Value = 2.7583717277179738
